# Se quitó COB Acapulco
# Removing the Acapulco "COB" adjusts ACAPULCO Ventas (col D) downward for the
# affected weeks, and redistributes the corresponding Tours (col E) total for
# week "2019 - 23" across the other plazas.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D77").Value = 25
$ws.Range("D79").Value = 29
$ws.Range("D80").Value = 48
$ws.Range("D82").Value = 80
$ws.Range("D83").Value = 71
$ws.Range("D84").Value = 87
$ws.Range("D85").Value = 96
$ws.Range("D87").Value = 43
$ws.Range("D89").Value = 21
$ws.Range("D90").Value = 28
$ws.Range("D91").Value = 31
$ws.Range("D92").Value = 20
$ws.Range("D93").Value = 20
$ws.Range("D95").Value = 18
$ws.Range("D98").Value = 35
$ws.Range("D99").Value = 21
$ws.Range("D100").Value = 57
$ws.Range("D101").Value = 29
$ws.Range("D102").Value = 28
$ws.Range("D103").Value = 37
$ws.Range("D104").Value = 39
$ws.Range("D105").Value = 62
$ws.Range("D107").Value = 43
$ws.Range("D108").Value = 34
$ws.Range("D109").Value = 35
$ws.Range("D110").Value = 31
$ws.Range("D113").Value = 35
$ws.Range("D118").Value = 34
$ws.Range("D121").Value = 22
$ws.Range("D122").Value = 28
$ws.Range("D123").Value = 38
$ws.Range("D125").Value = 24
$ws.Range("D127").Value = 16
$ws.Range("D131").Value = 30
$ws.Range("E167").Value = 141
$ws.Range("E168").Value = 120
$ws.Range("E169").Value = 236
$ws.Range("E170").Value = 126
$ws.Range("E171").Value = 117
$ws.Range("E172").Value = 129
$ws.Range("D181").Value = 27
$ws.Range("E181").Value = 110
$ws.Range("D251").Value = 18
$ws.Range("D252").Value = 5
$ws.Range("D253").Value = 6
$ws.Range("D254").Value = 6
$ws.Range("D256").Value = 7
$ws.Range("D259").Value = 4
$ws.Range("D260").Value = 6
$ws.Range("D261").Value = 1
$ws.Range("D262").Value = 13
$ws.Range("D263").Value = 11
$ws.Range("D264").Value = 22
$ws.Range("D268").Value = 7
$ws.Range("D269").Value = 6
$ws.Range("D270").Value = 1
$ws.Range("D271").Value = 4
$ws.Range("D272").Value = 4
$ws.Range("D274").Value = 5
$ws.Range("D275").Value = 5
$ws.Range("D276").Value = 7
$ws.Range("D277").Value = 8
$ws.Range("D278").Value = 5
$ws.Range("D279").Value = 5
$ws.Range("D281").Value = 6
$ws.Range("D282").Value = 1
$ws.Range("D283").Value = 8
$ws.Range("D284").Value = 12
$ws.Range("D285").Value = 27
$ws.Range("D286").Value = 33
$ws.Range("D287").Value = 19
$ws.Range("D288").Value = 5
$ws.Range("D289").Value = 3
$ws.Range("D290").Value = 4
$ws.Range("D291").Value = 15
$ws.Range("D292").Value = 6
$ws.Range("D293").Value = 4
$ws.Range("D294").Value = 8
$ws.Range("D295").Value = 7
$ws.Range("D296").Value = 12
$ws.Range("D297").Value = 9
$ws.Range("D298").Value = 16
$ws.Range("D299").Value = 28
$ws.Range("D301").Value = 4
$ws.Range("D302").Value = 5
$ws.Range("D303").Value = 10
$ws.Range("D304").Value = 7
$ws.Range("D305").Value = 5
$ws.Range("D306").Value = 9
$ws.Range("D307").Value = 7
$ws.Range("D308").Value = 6
$ws.Range("D309").Value = 5
$ws.Range("D310").Value = 8
$ws.Range("D311").Value = 12
$ws.Range("D313").Value = 10
$ws.Range("D314").Value = 13
$ws.Range("D315").Value = 16
$ws.Range("D316").Value = 21
$ws.Range("D317").Value = 17
$ws.Range("D318").Value = 9
$ws.Range("D319").Value = 3
$ws.Range("D320").Value = 8
$ws.Range("D321").Value = 7
$ws.Range("D322").Value = 5
$ws.Range("D323").Value = 8
$ws.Range("D324").Value = 5
$ws.Range("D325").Value = 4
$ws.Range("D326").Value = 6
$ws.Range("D327").Value = 3
$ws.Range("D328").Value = 5
$ws.Range("D329").Value = 8
$ws.Range("D330").Value = 8
$ws.Range("D331").Value = 9
$ws.Range("D333").Value = 5
$ws.Range("D336").Value = 14
$ws.Range("D337").Value = 33
$ws.Range("D339").Value = 17
$ws.Range("D341").Value = 8
$ws.Range("D342").Value = 12
$ws.Range("D348").Value = 13
$ws.Range("D350").Value = 10
$ws.Range("D352").Value = 14
$ws.Range("D353").Value = 24
$ws.Range("D354").Value = 38
$ws.Range("D355").Value = 9
$ws.Range("E361").Value = 45
$ws.Range("E508").Value = 8
$ws.Range("E685").Value = 29
$ws.Range("D860").Value = 8
$ws.Range("E860").Value = 54
$ws.Range("E979").Value = 19
$ws.Range("D1577").Value = 22
$ws.Range("E1577").Value = 130
$ws.Range("E1757").Value = 33
$ws.Range("D1937").Value = 38
$ws.Range("E1937").Value = 302
$ws.Range("E2111").Value = 56
